$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.863.40'
$ws.Range("E2").Value = '  -4.84%  '
$ws.Range("D3").Value = '2.213.83'
$ws.Range("E3").Value = '  -6.26%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.581'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.94%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.561'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0825'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.94%  '
$ws.Range("E14").Value = '  -3.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.862'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -11.57%  '
$ws.Range("D16").Value = '2.550.64'
$ws.Range("E16").Value = '  -6.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.08%  '
$ws.Range("D18").Value = '2.207.77'
$ws.Range("E18").Value = '  -6.78%  '
$ws.Range("D19").Value = '42.707.30'
$ws.Range("E19").Value = '  -5.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '0.0₃0961'
$ws.Range("E21").Value = '  -9.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -12.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '236.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.43%  '
$ws.Range("E26").Value = '  -7.54%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.29%  '
$ws.Range("E29").Value = '  -5.53%  '
$ws.Range("E30").Value = '  -13.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0871'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.91'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.68'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.75'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.18'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +15.20%  '
$ws.Range("E38").Value = '  -6.17%  '
$ws.Range("E39").Value = '  -6.23%  '
$ws.Range("E40").Value = '  -12.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0325'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.84%  '
$ws.Range("D43").Value = '1.873.64'
$ws.Range("E43").Value = '  +4.79%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.81%  '
$ws.Range("E47").Value = '  -9.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '60.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -12.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.62%  '
$ws.Range("E51").Value = '  -6.04%  '
